$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix hardcoded year in "HOMEROOM 2PINTAR 2022" -> "HOMEROOM 2PINTAR 2023"
$ws.Range("B8").Value = "HOMEROOM 2PINTAR 2023"

# Reset merit/demerit values that were leftover from previous period to 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("D36").Value = 0
